# Report diagrams and results
# Populate a side-by-side comparison block (columns T:Y) on Arkusz2 (sheet 2)
# mirroring the two existing result tables (B:D, rows 4-13 and rows 17-26),
# and link the new "\\" cells with external hyperlinks, matching the
# formatting/pattern already used on Arkusz1 (sheet 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# xlCenter
$xlCenter = -4108

# --- Block 1 (Euclidean/Places table, rows 4-13) -> mirrored into T4:Y13 ---
# T = this row's own k (column B)
# V = this block's own column C value (same row)
# X = the matching row's column C value from block 2 (row offset +13)
for ($r = 4; $r -le 13; $r++) {
    $k = $ws.Cells.Item($r, 2).Value2
    $r2 = $r + 13
    $v1 = $ws.Cells.Item($r, 3).Value2
    $v2 = $ws.Cells.Item($r2, 3).Value2

    $cT = $ws.Cells.Item($r, 20)
    $cT.HorizontalAlignment = $xlCenter
    $cT.Value = $k

    $cU = $ws.Cells.Item($r, 21)
    $cU.HorizontalAlignment = $xlCenter
    $cU.NumberFormat = "0.00%"
    $cU.Value = "&"

    $cV = $ws.Cells.Item($r, 22)
    $cV.HorizontalAlignment = $xlCenter
    $cV.NumberFormat = "0.00%"
    $cV.Value = $v1

    $cW = $ws.Cells.Item($r, 23)
    $cW.HorizontalAlignment = $xlCenter
    $cW.NumberFormat = "0.00%"
    $cW.Value = "&"

    $cX = $ws.Cells.Item($r, 24)
    $cX.HorizontalAlignment = $xlCenter
    $cX.NumberFormat = "0.00%"
    $cX.Value = $v2

    $cY = $ws.Cells.Item($r, 25)
    $cY.Value = "\\"
}

# --- Block 2 (Topics table, rows 17-26) -> mirrored into T17:Y26 ---
# T = this row's own k (column B)
# V = the matching row's column D value from block 1 (row offset -13)
# X = this block's own column D value (same row)
for ($r = 17; $r -le 26; $r++) {
    $k = $ws.Cells.Item($r, 2).Value2
    $r1 = $r - 13
    $v1 = $ws.Cells.Item($r1, 4).Value2
    $v2 = $ws.Cells.Item($r, 4).Value2

    $cT = $ws.Cells.Item($r, 20)
    $cT.HorizontalAlignment = $xlCenter
    $cT.Value = $k

    $cU = $ws.Cells.Item($r, 21)
    $cU.HorizontalAlignment = $xlCenter
    $cU.NumberFormat = "0.00%"
    $cU.Value = "&"

    $cV = $ws.Cells.Item($r, 22)
    $cV.HorizontalAlignment = $xlCenter
    $cV.NumberFormat = "0.00%"
    $cV.Value = $v1

    $cW = $ws.Cells.Item($r, 23)
    $cW.HorizontalAlignment = $xlCenter
    $cW.NumberFormat = "0.00%"
    $cW.Value = "&"

    $cX = $ws.Cells.Item($r, 24)
    $cX.HorizontalAlignment = $xlCenter
    $cX.NumberFormat = "0.00%"
    $cX.Value = $v2

    $cY = $ws.Cells.Item($r, 25)
    $cY.Value = "\\"
}

# Hyperlinks over the "\\" column, grouped the same way Arkusz1 groups its
# analogous hyperlink cells: a single-cell hyperlink on the first row of each
# block, then one range hyperlink covering the rest of that block.
$hy4 = $ws.Range("Y4")
$ws.Hyperlinks.Add($hy4, "\")
$hy4.Style = "Hiperłącze"

$hy5_13 = $ws.Range("Y5:Y13")
$ws.Hyperlinks.Add($hy5_13, "\", "", "", "\\")
$hy5_13.Style = "Hiperłącze"

$hy17 = $ws.Range("Y17")
$ws.Hyperlinks.Add($hy17, "\")
$hy17.Style = "Hiperłącze"

$hy18_26 = $ws.Range("Y18:Y26")
$ws.Hyperlinks.Add($hy18_26, "\", "", "", "\\")
$hy18_26.Style = "Hiperłącze"

# Scroll / selection, matching the saved view state (topLeftCell ~ column C,
# active cell T17 with the new block selected).
$ws.Activate()
$ws.Range("T17:Y26").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
